$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 11.59547033333333
$ws.Range("H2").Value = 34.786411
$ws.Range("I2").Value = 0.2556603672772884
$ws.Range("J2").Value = 0.2556603672772884
$ws.Range("M2").Value = 38.745275
$ws.Range("N2").Value = 116.235825
$ws.Range("O2").Value = 0.3160319337595895
$ws.Range("P2").Value = 0.3160319337595895
$ws.Range("Q2").Value = 449.2696868193416
$ws.Range("R2").Value = 4043.427181374075
$ws.Range("S2").Value = 0.08079684025632833
$ws.Range("T2").Value = 0.08079684025632831
$ws.Range("G3").Value = 11.59547033333333
$ws.Range("H3").Value = 34.786411
$ws.Range("I3").Value = 0.2556603672772884
$ws.Range("J3").Value = 0.2556603672772884
$ws.Range("O3").Value = 0.4383510712400526
$ws.Range("P3").Value = 0.4383510712400526
$ws.Range("Q3").Value = 623.1580655477528
$ws.Range("R3").Value = 5608.422589929775
$ws.Range("S3").Value = 0.1120689958696247
$ws.Range("T3").Value = 0.1120689958696246
$ws.Range("G4").Value = 11.59547033333333
$ws.Range("H4").Value = 34.786411
$ws.Range("I4").Value = 0.2556603672772884
$ws.Range("J4").Value = 0.2556603672772884
$ws.Range("M4").Value = 17.38482166666667
$ws.Range("N4").Value = 52.154465
$ws.Range("O4").Value = 0.1418020341675798
$ws.Range("P4").Value = 0.1418020341675798
$ws.Range("Q4").Value = 201.5851838861239
$ws.Range("R4").Value = 1814.266654975115
$ws.Range("S4").Value = 0.03625316013595005
$ws.Range("T4").Value = 0.03625316013595005
$ws.Range("G5").Value = 11.59547033333333
$ws.Range("H5").Value = 34.786411
$ws.Range("I5").Value = 0.2556603672772884
$ws.Range("J5").Value = 0.2556603672772884
$ws.Range("M5").Value = 12.72763533333333
$ws.Range("N5").Value = 38.182906
$ws.Range("O5").Value = 0.103814960832778
$ws.Range("P5").Value = 0.103814960832778
$ws.Range("Q5").Value = 147.5829179211518
$ws.Range("R5").Value = 1328.246261290366
$ws.Range("S5").Value = 0.02654137101538532
$ws.Range("T5").Value = 0.02654137101538532
$ws.Range("I6").Value = 0.2896792350628708
$ws.Range("J6").Value = 0.2896792350628708
$ws.Range("M6").Value = 38.745275
$ws.Range("N6").Value = 116.235825
$ws.Range("O6").Value = 0.3160319337595895
$ws.Range("P6").Value = 0.3160319337595895
$ws.Range("Q6").Value = 509.0507402487166
$ws.Range("R6").Value = 4581.456662238449
$ws.Range("S6").Value = 0.09154788882691775
$ws.Range("T6").Value = 0.09154788882691774
$ws.Range("I7").Value = 0.2896792350628708
$ws.Range("J7").Value = 0.2896792350628708
$ws.Range("O7").Value = 0.4383510712400526
$ws.Range("P7").Value = 0.4383510712400526
$ws.Range("S7").Value = 0.1269812030058084
$ws.Range("T7").Value = 0.1269812030058084
$ws.Range("I8").Value = 0.2896792350628708
$ws.Range("J8").Value = 0.2896792350628708
$ws.Range("M8").Value = 17.38482166666667
$ws.Range("N8").Value = 52.154465
$ws.Range("O8").Value = 0.1418020341675798
$ws.Range("P8").Value = 0.1418020341675798
$ws.Range("Q8").Value = 228.4086598561656
$ws.Range("R8").Value = 2055.67793870549
$ws.Range("S8").Value = 0.04107710478802359
$ws.Range("T8").Value = 0.04107710478802359
$ws.Range("I9").Value = 0.2896792350628708
$ws.Range("J9").Value = 0.2896792350628708
$ws.Range("M9").Value = 12.72763533333333
$ws.Range("N9").Value = 38.182906
$ws.Range("O9").Value = 0.103814960832778
$ws.Range("P9").Value = 0.103814960832778
$ws.Range("Q9").Value = 167.2207046678351
$ws.Range("R9").Value = 1504.986342010516
$ws.Range("S9").Value = 0.03007303844212101
$ws.Range("T9").Value = 0.03007303844212101
$ws.Range("G10").Value = 7.691597000000001
$ws.Range("H10").Value = 23.074791
$ws.Range("I10").Value = 0.1695866107574785
$ws.Range("J10").Value = 0.1695866107574785
$ws.Range("M10").Value = 38.745275
$ws.Range("N10").Value = 116.235825
$ws.Range("O10").Value = 0.3160319337595895
$ws.Range("P10").Value = 0.3160319337595895
$ws.Range("Q10").Value = 298.013040954175
$ws.Range("R10").Value = 2682.117368587575
$ws.Range("S10").Value = 0.05359478453742074
$ws.Range("T10").Value = 0.05359478453742073
$ws.Range("G11").Value = 7.691597000000001
$ws.Range("H11").Value = 23.074791
$ws.Range("I11").Value = 0.1695866107574785
$ws.Range("J11").Value = 0.1695866107574785
$ws.Range("O11").Value = 0.4383510712400526
$ws.Range("P11").Value = 0.4383510712400526
$ws.Range("Q11").Value = 413.3580242721417
$ws.Range("R11").Value = 3720.222218449275
$ws.Range("S11").Value = 0.07433847249351053
$ws.Range("T11").Value = 0.07433847249351053
$ws.Range("G12").Value = 7.691597000000001
$ws.Range("H12").Value = 23.074791
$ws.Range("I12").Value = 0.1695866107574785
$ws.Range("J12").Value = 0.1695866107574785
$ws.Range("M12").Value = 17.38482166666667
$ws.Range("N12").Value = 52.154465
$ws.Range("O12").Value = 0.1418020341675798
$ws.Range("P12").Value = 0.1418020341675798
$ws.Range("Q12").Value = 133.7170421768683
$ws.Range("R12").Value = 1203.453379591815
$ws.Range("S12").Value = 0.02404772637299603
$ws.Range("T12").Value = 0.02404772637299603
$ws.Range("G13").Value = 7.691597000000001
$ws.Range("H13").Value = 23.074791
$ws.Range("I13").Value = 0.1695866107574785
$ws.Range("J13").Value = 0.1695866107574785
$ws.Range("M13").Value = 12.72763533333333
$ws.Range("N13").Value = 38.182906
$ws.Range("O13").Value = 0.103814960832778
$ws.Range("P13").Value = 0.103814960832778
$ws.Range("Q13").Value = 97.89584174696068
$ws.Range("R13").Value = 881.0625757226461
$ws.Range("S13").Value = 0.01760562735355119
$ws.Range("T13").Value = 0.01760562735355119
$ws.Range("G14").Value = 12.92951533333333
$ws.Range("H14").Value = 38.788546
$ws.Range("I14").Value = 0.2850737869023623
$ws.Range("J14").Value = 0.2850737869023623
$ws.Range("M14").Value = 38.745275
$ws.Range("N14").Value = 116.235825
$ws.Range("O14").Value = 0.3160319337595895
$ws.Range("P14").Value = 0.3160319337595895
$ws.Range("Q14").Value = 500.9576272067167
$ws.Range("R14").Value = 4508.618644860449
$ws.Range("S14").Value = 0.09009242013892273
$ws.Range("T14").Value = 0.09009242013892271
$ws.Range("G15").Value = 12.92951533333333
$ws.Range("H15").Value = 38.788546
$ws.Range("I15").Value = 0.2850737869023623
$ws.Range("J15").Value = 0.2850737869023623
$ws.Range("O15").Value = 0.4383510712400526
$ws.Range("P15").Value = 0.4383510712400526
$ws.Range("Q15").Value = 694.8516560322945
$ws.Range("R15").Value = 6253.664904290649
$ws.Range("S15").Value = 0.124962399871109
$ws.Range("T15").Value = 0.124962399871109
$ws.Range("G16").Value = 12.92951533333333
$ws.Range("H16").Value = 38.788546
$ws.Range("I16").Value = 0.2850737869023623
$ws.Range("J16").Value = 0.2850737869023623
$ws.Range("M16").Value = 17.38482166666667
$ws.Range("N16").Value = 52.154465
$ws.Range("O16").Value = 0.1418020341675798
$ws.Range("P16").Value = 0.1418020341675798
$ws.Range("Q16").Value = 224.7773183064322
$ws.Range("R16").Value = 2022.99586475789
$ws.Range("S16").Value = 0.04042404287061016
$ws.Range("T16").Value = 0.04042404287061016
$ws.Range("G17").Value = 12.92951533333333
$ws.Range("H17").Value = 38.788546
$ws.Range("I17").Value = 0.2850737869023623
$ws.Range("J17").Value = 0.2850737869023623
$ws.Range("M17").Value = 12.72763533333333
$ws.Range("N17").Value = 38.182906
$ws.Range("O17").Value = 0.103814960832778
$ws.Range("P17").Value = 0.103814960832778
$ws.Range("Q17").Value = 164.5621561994084
$ws.Range("R17").Value = 1481.059405794676
$ws.Range("S17").Value = 0.02959492402172044
$ws.Range("T17").Value = 0.02959492402172044
